# Refresh the crypto symbol/price listing (GitHub Actions scheduled update).
#
# The sheet stores every data cell as literal text (prices such as "246.15"
# are strings, not numbers), so for any new value that *looks* numeric we
# must format the cell as Text ("@") before assigning it - otherwise Excel
# would silently coerce it to a real number and drop things like trailing
# zeros (e.g. "0.05730" -> 0.0573).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
}

function Set-StringValue($addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- Price refresh for the top-of-list coins (rows 2-8) ---
Set-TextValue "D2" "246.33"
Set-TextValue "D3" "22.76"
Set-TextValue "D4" "5.276"
Set-TextValue "D5" "0.05730"
Set-TextValue "D6" "3.436"
Set-TextValue "D7" "0.8094"
Set-TextValue "D8" "0.8818"

# --- Rows 9-11: ranking reshuffled, new coin "One" enters at rank 8,
#     pushing WazirX / MandalaExchangeToken / ProBitToken down one slot ---
Set-StringValue "B9" "One"
Set-StringValue "C9" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue   "D9" "0.01097"
Set-StringValue "E9" "8OneONEBestin24h"

Set-StringValue "B10" "WazirX"
Set-StringValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue   "D10" "0.1427"
Set-StringValue "E10" "9WazirXWRX"

Set-StringValue "B11" "MandalaExchangeToken"
Set-StringValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue   "D11" "0.07366"
Set-StringValue "E11" "10MandalaExchangeTokenMDX"

# --- Rows 12-17: price refresh only, ranking unchanged ---
Set-TextValue "D12" "0.03035"
Set-TextValue "D13" "0.03116"
Set-TextValue "D14" "0.09383"
Set-TextValue "D15" "3.927"
Set-TextValue "D16" "0.001588"
Set-TextValue "D17" "0.04798"

# --- Rows 18-26: ranking reshuffled, "One" drops out of this block and
#     every other coin shifts up one slot, with ProBitToken re-entering
#     at rank 25 ---
Set-StringValue "B18" "TigerCash"
Set-StringValue "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue   "D18" "0.006140"
Set-StringValue "E18" "17TigerCashTCH"

Set-StringValue "B19" "HotbitToken"
Set-StringValue "C19" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue   "D19" "0.005113"
Set-StringValue "E19" "18HotbitTokenHTB"

Set-StringValue "B20" "BitKan"
Set-StringValue "C20" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue   "D20" "0.0009963"
Set-StringValue "E20" "19BitKanKAN"

Set-StringValue "B21" "NitroEx"
Set-StringValue "C21" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue   "D21" "0.0001500"
Set-StringValue "E21" "20NitroExNTX"

Set-StringValue "B22" "LEO"
Set-StringValue "C22" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue   "D22" "3.745"
Set-StringValue "E22" "21LEOLEO"

Set-StringValue "B23" "KuCoinToken"
Set-StringValue "C23" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue   "D23" "6.315"
Set-StringValue "E23" "22KuCoinTokenKCS"

Set-StringValue "B24" "BTSEToken"
Set-StringValue "C24" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue   "D24" "2.186"
Set-StringValue "E24" "23BTSETokenBTSE"

Set-StringValue "B25" "BitpandaEcosystemToken"
Set-StringValue "C25" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue   "D25" "0.3280"
Set-StringValue "E25" "24BitpandaEcosystemTokenBEST"

Set-StringValue "B26" "ProBitToken"
Set-StringValue "C26" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue   "D26" "0.1350"
Set-StringValue "E26" "25ProBitTokenPROB"

# --- Rows 40-44: price refresh only ---
Set-TextValue "D40" "0.03911"
Set-TextValue "D41" "0.006713"
Set-TextValue "D42" "0.1066"
Set-TextValue "D43" "0.003201"
Set-TextValue "D44" "0.007488"

# --- Row 47: "best in 24h" badge removed from CoinbaseStockToken's volume label ---
Set-StringValue "E47" "46CoinbaseStockTokenCOIN"

# --- Row 48: price refresh + "worst in 24h" badge added to BOLO's volume label ---
Set-TextValue   "D48" "0.1749"
Set-StringValue "E48" "47BOLOBOLOWorstin24h"
